# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.951.11'
$ws.Range("E2").Value = '  +2.80%  '

$ws.Range("D3").Value = '2.608.40'
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''572.43'
$ws.Range("E5").Value = '  -0.21%  '

$ws.Range("D6").Value = '''143.59'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").Value = '''0.996'
$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("E8").Value = '  +1.02%  '

$ws.Range("D9").Value = '2.633.59'
$ws.Range("E9").Value = '  +1.76%  '

$ws.Range("D10").Value = '''6.53'
$ws.Range("E10").Value = '  -2.57%  '

$ws.Range("E11").Value = '  +3.01%  '

$ws.Range("E12").Value = '  -3.31%  '

$ws.Range("D13").Value = '''0.368'
$ws.Range("E13").Value = '  +6.41%  '

$ws.Range("D14").Value = '3.086.69'

$ws.Range("D15").Value = '60.946.51'
$ws.Range("E15").Value = '  +2.79%  '

$ws.Range("D16").Value = '''23.54'
$ws.Range("E16").Value = '  +4.60%  '

$ws.Range("E17").Value = '  +2.73%  '

$ws.Range("D18").Value = '2.631.84'
$ws.Range("E18").Value = '  +1.85%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''4.67'
$ws.Range("E19").Value = '  +3.11%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''11.21'
$ws.Range("E20").Value = '  +9.13%  '

$ws.Range("D21").Value = '''352.31'
$ws.Range("E21").Value = '  +3.72%  '

$ws.Range("D22").Value = '''7.18'
$ws.Range("E22").Value = '  +14.51%  '

$ws.Range("E23").Value = '  +0.34%  '

$ws.Range("D24").Value = '''0.520'
$ws.Range("E24").Value = '  +13.16%  '

$ws.Range("D25").Value = '''64.15'
$ws.Range("E25").Value = '  -0.79%  '

$ws.Range("D26").Value = '''0.998'
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").Value = '''0.161'
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("E28").Value = '  +6.14%  '

$ws.Range("E29").Value = '  +1.63%  '

$ws.Range("E30").Value = '  +6.70%  '

$ws.Range("D31").Value = '''0.997'
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").Value = '''6.34'
$ws.Range("E32").Value = '  +4.61%  '

$ws.Range("D33").Value = '''160.89'
$ws.Range("E33").Value = '  +1.75%  '

$ws.Range("D34").Value = '''19.56'
$ws.Range("E34").Value = '  +2.76%  '

$ws.Range("D35").Value = '''4.28'
$ws.Range("E35").Value = '  +5.90%  '

$ws.Range("D36").Value = '''0.960'
$ws.Range("E36").Value = '  +9.29%  '

$ws.Range("E37").Value = '  +4.12%  '

$ws.Range("E38").Value = '  +6.38%  '

$ws.Range("D39").Value = '''37.78'
$ws.Range("E39").Value = '  +1.44%  '

$ws.Range("E40").Value = '  -2.11%  '

$ws.Range("D41").Value = '''3.81'
$ws.Range("E41").Value = '  +3.63%  '

$ws.Range("D42").Value = '''298.70'
$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("D43").Value = '''140.55'
$ws.Range("E43").Value = '  +9.49%  '

$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("D45").Value = '''0.995'
$ws.Range("E45").Value = '  -0.27%  '

$ws.Range("E46").Value = '  +2.17%  '

$ws.Range("D47").Value = '''0.0551'
$ws.Range("E47").Value = '  +2.60%  '

$ws.Range("D48").Value = '''0.0241'
$ws.Range("E48").Value = '  +3.45%  '

$ws.Range("D49").Value = '''10.71'
$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("D50").Value = '''19.65'
$ws.Range("E50").Value = '  +5.77%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '''4.86'
$ws.Range("E51").Value = '  +7.78%  '
